$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("O34").Value = "Green"
